$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B65").Value = "22f90730e0e3a4833b18ff80dfc6da2c"
$ws.Range("B80").Value = "5c7461dca3df71bee93b2ffb4c7aff6b"
$ws.Range("B81").Value = "49281ee1136c065ad894839c40f98be3"
$ws.Range("B113").Value = "d29b6e376f4ab6820f3ee6102e491c52"
$ws.Range("B117").Value = "81a54580528beaa3746c0be2eb8a639f"
$ws.Range("B122").Value = "ee278e6bca7a8de6505b4498ce294b23"
$ws.Range("B163").Value = "3c645c889328825fd153ef3d644e8721"
$ws.Range("B187").Value = "bbe42b101b0df60ce51714a09105540e"
$ws.Range("B527").Value = "bbe42b101b0df60ce51714a09105540e"
$ws.Range("B342").Value = "b105fc2bcbf5ba38e657e44c9d60c42d"
$ws.Range("B404").Value = "1c944e1273c1f4f89c6b614f6f639dcd"
$ws.Range("B487").Value = "ebc4ac9707a429ede3a6c3dd3f2f2a4f"
$ws.Range("B507").Value = "af81f829b84b7a5c9ca178ceb0f4a0f8"
$ws.Range("B548").Value = "d02747c32d0174eecaf3932cfc67d53a"
$ws.Range("B574").Value = "f3ec5be5d08573163925e4f336c4139c"
$ws.Range("B619").Value = "e3ee95ef384d09352f2806899d18ac19"
$ws.Range("B697").Value = "3bfb87aceee6eb9ba861adc6a9cb0d3f"
$ws.Range("B700").Value = "82b72890bcd18c3586e86b913781f104"
$ws.Range("B763").Value = "e66db9edeb85723f367334b05a32f91c"
$ws.Range("B822").Value = "27ef180bf5c47eebd4d6b6059e4f00b3"
$ws.Range("B890").Value = "6d01a0f675fe2d3fc4c8a159403d3cf2"
$ws.Range("B946").Value = "68d93f49edbf51e7b78c496557f4b6ba"
